$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8 and 20 Trigger column ("N" -> "Y")
$ws.Range("C8").Value = "Y"
$ws.Range("C20").Value = "Y"

# Add hyperlinks for the verifyURL rows (Extent Report login links).
# F28 keeps its displayed text ("https://github.com/login"); F22 and F10
# get both their target and their displayed text updated to the new URL.
$ws.Hyperlinks.Add($ws.Range("F28"), "https://github.com/login") | Out-Null

$ws.Range("F22").Value = "https://github.com/login04"
$ws.Hyperlinks.Add($ws.Range("F22"), "https://github.com/login04") | Out-Null

$ws.Range("F10").Value = "https://github.com/login04"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://github.com/login04") | Out-Null

# Move the saved selection to H4
$ws.Range("H4").Select() | Out-Null

Write-Output "done"
